# Working on MVP and adding devcontainer
#
# 1. slide2 ("Test Text" -> "Test Text Content") + new "Test some text"
#    content placeholder.
# 2. New slide3: "Bulleted List Content" title, with a two-line bulleted
#    body ("some important info" / "important sub-bullet" at level 2).
# 3. New slide4: "Picture Content" title (Title Only layout).

$p = $ppt.ActivePresentation

# --- Slide 2: extend the title and fill in the content placeholder ---
$slide2 = $p.Slides.Item(2)

# Append " Content" onto the existing "Test Text" run rather than
# overwriting the whole range, so we don't disturb the existing run.
[void]$slide2.Shapes.Item(1).TextFrame.TextRange.InsertAfter(" Content")

# The second placeholder ("Content Placeholder 2") starts out empty;
# insert its text directly.
[void]$slide2.Shapes.Item(2).TextFrame.TextRange.InsertAfter("Test some text")

# --- Slide 3 (new): "Bulleted List Content" with a sub-bullet ---
$slide3 = $p.Slides.Add(3, 2)   # ppLayoutText -> "Title and Content"

[void]$slide3.Shapes.Item(1).TextFrame.TextRange.InsertAfter("Bulleted List Content")

$body3 = $slide3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "some important info`rimportant sub-bullet"
$body3.Paragraphs(2).IndentLevel = 2

# --- Slide 4 (new): "Picture Content" title only ---
$slide4 = $p.Slides.Add(4, 6)   # ppLayoutTitleOnly -> "Title Only"

[void]$slide4.Shapes.Item(1).TextFrame.TextRange.InsertAfter("Picture Content")
